# Regenerate save_data: column G ("K") values recalculated (K instead of Strike#),
# std/mean regenerated and s_vals written. Apply the updated K values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 3
    12 = 1
    13 = 0
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
